# Add Board dinner receipt
# Account 6860 "Møte, kurs, oppdatering" (row 67): the actual 2021 amount (column E)
# increases by 3000 kr to reflect a newly recorded Board dinner receipt.
# All the other values on the sheet (F67, E78/F78, E80/F80, E82/F82) are formulas
# that recalculate automatically from this single change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E67").Value = 5509.93

# Restore the last-used selection/active cell as recorded in the saved file.
$ws.Range("G58").Select()
